## Generate Report for Archive
##
## 1) The status text "Ready for handoff" becomes "In Translation" everywhere
##    it appears (Overview!E2, Overview!F2, zh-cn!C2, de-de!C2 all shared the
##    same string).
## 2) The "Latest Handoff Datetime" / status columns are narrower now:
##      - Overview sheet: columns E and F (zh-cn / de-de status columns)
##      - zh-cn sheet: column C (Status)
##      - de-de sheet: column C (Status)

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Update status text (shared across sheets) ---
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- Narrow the status columns ---
$wsOverview.Columns.Item(5).ColumnWidth = 13.4101845877511
$wsOverview.Columns.Item(6).ColumnWidth = 13.4101845877511

$wsZhCn.Columns.Item(3).ColumnWidth = 13.4101845877511

$wsDeDe.Columns.Item(3).ColumnWidth = 13.4101845877511
